$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at 56, shifting the existing rows 56+ down by one.
$ws.Rows("56:56").Insert(-4121)  # xlShiftDown
$ws.Range("A56:P56").Clear()

# Populate the new "Wirksworth" pub run entry in row 56.
$ws.Range("A56").Value = 43845
$ws.Range("B56").Value = "The Malt Shovel"
$ws.Range("C56").Value = "Wirksworth Moor"
$ws.Range("D56").Value = "start/end at pub"
$ws.Range("E56").Value = 2.71
$ws.Range("F56").Value = 0.0408912037037037
$ws.Range("G56").Formula = "=F56/E56"
$ws.Range("H56").Value = 1
$ws.Range("J56").Value = 1
$ws.Range("N56").Value = 1
$ws.Range("O56").Value = "Astronomical observations, Black Rocks"
$ws.Range("P56").Formula = "=SUM(H56:N56)*E56"

# Match formatting of the row above (the previous last entry) for the new row.
$ws.Range("A55:P55").Copy()
$ws.Range("A56:P56").PasteSpecial(-4122)  # xlPasteFormats

# Columns that should stay blank for this entry (the row-55 template had
# values there that don't apply here) - remove the stray empty styled cells.
$ws.Range("I56").Clear()
$ws.Range("K56").Clear()
$ws.Range("L56").Clear()
$ws.Range("M56").Clear()

$ws.Rows("56:56").RowHeight = 11.25

# Update the selected cell to reflect where the new entry was added.
$ws.Range("A56").Select()
